# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new labels in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold/border/centered style).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows (2-64): same Wins/Losses/Ties record for every player.
for ($r = 2; $r -le 64; $r++) {
    $ws.Cells.Item($r, 30).Value = 50
    $ws.Cells.Item($r, 31).Value = 112
    $ws.Cells.Item($r, 32).Value = 0
}
